$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.185301542282104
$ws.Range("B1").Value = 4.095221996307373
$ws.Range("C1").Value = 3.100350618362427
$ws.Range("D1").Value = 2.940971851348877
$ws.Range("E1").Value = 2.398838043212891
